$d = $word.ActiveDocument

$replacements = @(
    @("831÷5=", "159÷8="),
    @("555÷4=", "793÷7="),
    @("230÷9=", "732÷8="),
    @("461÷3=", "857÷2="),
    @("551÷4=", "323÷7="),
    @("547÷6=", "949÷7="),
    @("942÷7=", "324÷6="),
    @("816÷6=", "404÷8="),
    @("576÷8=", "240÷9="),
    @("181÷3=", "286÷7="),
    @("788÷9=", "926÷9="),
    @("381÷2=", "345÷2="),
    @("719÷2=", "993÷7="),
    @("416÷4=", "421÷3="),
    @("708÷9=", "566÷4="),
    @("161÷7=", "729÷2="),
    @("535÷3=", "187÷8="),
    @("943÷9=", "458÷8="),
    @("159÷7=", "958÷7="),
    @("443÷8=", "536÷8="),
    @("838÷8=", "575÷5="),
    @("794÷9=", "527÷3="),
    @("698÷2=", "470÷3="),
    @("725÷3=", "949÷4="),
    @("153÷7=", "311÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
